$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting for numeric-looking values
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '61.030.27'
$ws.Range("E2").Value = '  -1.89%  '

# Row 3
$ws.Range("D3").Value = '2.433.36'
$ws.Range("E3").Value = '  -0.39%  '

# Row 4
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").Value = '572.80'
$ws.Range("E5").Value = '  -1.55%  '

# Row 6
$ws.Range("D6").Value = '140.65'
$ws.Range("E6").Value = '  -1.85%  '

# Row 7
$ws.Range("E7").Value = '  +0.19%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("D9").Value = '2.419.85'
$ws.Range("E9").Value = '  -0.85%  '

# Row 10
$ws.Range("E10").Value = '  +1.55%  '

# Row 11
$ws.Range("E11").Value = '  +0.97%  '

# Row 12
$ws.Range("D12").Value = '5.13'
$ws.Range("E12").Value = '  -1.45%  '

# Row 13
$ws.Range("E13").Value = '  -1.32%  '

# Row 14
$ws.Range("E14").Value = '  -1.17%  '

# Row 15
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '0.0000171'
$ws.Range("E15").Value = '  -0.90%  '

# Row 16
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '2.860.40'

# Row 17
$ws.Range("D17").Value = '60.945.47'
$ws.Range("E17").Value = '  -1.84%  '

# Row 18
$ws.Range("D18").Value = '2.413.92'
$ws.Range("E18").Value = '  -0.84%  '

# Row 19
$ws.Range("E19").Value = '  -2.89%  '

# Row 20
$ws.Range("D20").Value = '7.28'
$ws.Range("E20").Value = '  +2.52%  '

# Row 21
$ws.Range("D21").Value = '323.65'
$ws.Range("E21").Value = '  -2.07%  '

# Row 22
$ws.Range("D22").Value = '4.04'
$ws.Range("E22").Value = '  -1.58%  '

# Row 23
$ws.Range("D23").Value = '6.14'
$ws.Range("E23").Value = '  +2.81%  '

# Row 24
$ws.Range("E24").Value = '  +0.17%  '

# Row 25
$ws.Range("E25").Value = '  -3.31%  '

# Row 26
$ws.Range("D26").Value = '65.22'
$ws.Range("E26").Value = '  -0.90%  '

# Row 27
$ws.Range("D27").Value = '8.89'
$ws.Range("E27").Value = '  -5.02%  '

# Row 28
$ws.Range("D28").Value = '576.03'
$ws.Range("E28").Value = '  -6.88%  '

# Row 29
$ws.Range("E29").Value = '  -0.29%  '

# Row 30
$ws.Range("E30").Value = '  -0.06%  '

# Row 31
$ws.Range("D31").Value = '0.0₃0916'

# Row 32
$ws.Range("D32").Value = '7.89'
$ws.Range("E32").Value = '  -1.30%  '

# Row 33
$ws.Range("D33").Value = '1.35'
$ws.Range("E33").Value = '  -5.33%  '

# Row 34
$ws.Range("E34").Value = '  -1.97%  '

# Row 35
$ws.Range("E35").Value = '  -6.22%  '

# Row 36
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.18%  '

# Row 37
$ws.Range("D37").Value = '4.62'
$ws.Range("E37").Value = '  -5.81%  '

# Row 38
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").Value = '0.369'
$ws.Range("E38").Value = '  -1.74%  '

# Row 39
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = '151.20'
$ws.Range("E39").Value = '  +0.01%  '

# Row 40
$ws.Range("E40").Value = '  -3.53%  '

# Row 41
$ws.Range("D41").Value = '18.27'
$ws.Range("E41").Value = '  -0.31%  '

# Row 42
$ws.Range("D42").Value = '5.11'
$ws.Range("E42").Value = '  -2.29%  '

# Row 43
$ws.Range("E43").Value = '  +0.03%  '

# Row 44
$ws.Range("E44").Value = '  -1.86%  '

# Row 45
$ws.Range("D45").Value = '1.66'
$ws.Range("E45").Value = '  -5.94%  '

# Row 46
$ws.Range("D46").Value = '2.35'
$ws.Range("E46").Value = '  -4.50%  '

# Row 47
$ws.Range("D47").Value = '0.0₆0285'
$ws.Range("E47").Value = '  +25.55%  '

# Row 48
$ws.Range("D48").Value = '141.31'
$ws.Range("E48").Value = '  -1.24%  '

# Row 49
$ws.Range("D49").Value = '3.53'
$ws.Range("E49").Value = '  -2.56%  '

# Row 50
$ws.Range("E50").Value = '  -0.67%  '

# Row 51
$ws.Range("D51").Value = '19.57'
$ws.Range("E51").Value = '  +0.44%  '
